# Update the per-innings batting stats (runs, balls, fours, sixes) for
# Kieron Pollard. The underlying data rows were re-ordered/re-shuffled;
# this writes the resulting values directly into C2:F11.
#
# The source sheet stores these numbers as text (number-stored-as-text,
# see the <ignoredErrors> entry covering A1:F11), so force the number
# format to Text ("@") before assigning the values - otherwise Excel
# would happily reinterpret "9" etc. as a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("9", "4", "2", "0")
    3  = @("0", "2", "0", "0")
    4  = @("34", "12", "1", "4")
    5  = @("60", "24", "3", "5")
    6  = @("41", "25", "2", "4")
    7  = @("11", "14", "1", "0")
    8  = @("13", "7", "1", "0")
    9  = @("47", "20", "3", "4")
    10 = @("18", "14", "1", "1")
    11 = @("25", "13", "0", "3")
}

$columns = @("C", "D", "E", "F")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cell = $ws.Range("$($columns[$i])$row")
        $cell.NumberFormat = "@"
        $cell.Value = $values[$i]
    }
}
